# Update the date heading, then the 25 division-problem answer cells in
# the single table. Cells are addressed by (row, column) position rather
# than by old-text matching, since some new values collide with other
# cells' old values (e.g. "972÷5=194, 2" is both an old value in one cell
# and the new value of a different cell) — positional addressing avoids
# any ambiguity that a text-based Find/Replace could introduce.

$d = $word.ActiveDocument

# Heading date line (first paragraph).
$d.Paragraphs.Item(1).Range.Text = "2024-11-22 Friday"

$t = $d.Tables.Item(1)

# Row 1 (table header row of answers)
$t.Cell(1, 1).Range.Text = "920÷4=230, 0"
$t.Cell(1, 2).Range.Text = "288÷3=96, 0"
$t.Cell(1, 3).Range.Text = "147÷9=16, 3"
$t.Cell(1, 4).Range.Text = "151÷2=75, 1"
$t.Cell(1, 5).Range.Text = "955÷2=477, 1"

# Row 5
$t.Cell(5, 1).Range.Text = "690÷7=98, 4"
$t.Cell(5, 2).Range.Text = "600÷4=150, 0"
$t.Cell(5, 3).Range.Text = "266÷4=66, 2"
$t.Cell(5, 4).Range.Text = "664÷3=221, 1"
$t.Cell(5, 5).Range.Text = "676÷9=75, 1"

# Row 9
$t.Cell(9, 1).Range.Text = "345÷2=172, 1"
$t.Cell(9, 2).Range.Text = "630÷7=90, 0"
$t.Cell(9, 3).Range.Text = "974÷6=162, 2"
$t.Cell(9, 4).Range.Text = "807÷8=100, 7"
$t.Cell(9, 5).Range.Text = "243÷3=81, 0"

# Row 13
$t.Cell(13, 1).Range.Text = "593÷6=98, 5"
$t.Cell(13, 2).Range.Text = "598÷8=74, 6"
$t.Cell(13, 3).Range.Text = "821÷8=102, 5"
$t.Cell(13, 4).Range.Text = "861÷3=287, 0"
$t.Cell(13, 5).Range.Text = "127÷7=18, 1"

# Row 17
$t.Cell(17, 1).Range.Text = "102÷6=17, 0"
$t.Cell(17, 2).Range.Text = "216÷9=24, 0"
$t.Cell(17, 3).Range.Text = "514÷5=102, 4"
$t.Cell(17, 4).Range.Text = "916÷6=152, 4"
$t.Cell(17, 5).Range.Text = "972÷5=194, 2"

Write-Host "Done updating date and 25 answer cells."
